$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to build each new value as a literal-text formula,
# then paste-special (values-only) into the target cell. This guarantees
# the result lands as text (matching the original inlineStr cells) instead
# of being auto-coerced into a number by a plain Range.Value assignment,
# and it leaves the target cell style untouched.
$scratch = $ws.Range("Z1")

$scratch.Formula = "=`"43.620.19`""
$scratch.Copy()
$ws.Range("D2").PasteSpecial(-4163)

$scratch.Formula = "=`"  +1.12%  `""
$scratch.Copy()
$ws.Range("E2").PasteSpecial(-4163)

$scratch.Formula = "=`"2.420.40`""
$scratch.Copy()
$ws.Range("D3").PasteSpecial(-4163)

$scratch.Formula = "=`"  +2.49%  `""
$scratch.Copy()
$ws.Range("E3").PasteSpecial(-4163)

$scratch.Formula = "=`"  +0.00%  `""
$scratch.Copy()
$ws.Range("E4").PasteSpecial(-4163)

$scratch.Formula = "=`"306.66`""
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)

$scratch.Formula = "=`"  +1.20%  `""
$scratch.Copy()
$ws.Range("E5").PasteSpecial(-4163)

$scratch.Formula = "=`"97.66`""
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)

$scratch.Formula = "=`"  +1.60%  `""
$scratch.Copy()
$ws.Range("E6").PasteSpecial(-4163)

$scratch.Formula = "=`"  +0.07%  `""
$scratch.Copy()
$ws.Range("E7").PasteSpecial(-4163)

$scratch.Formula = "=`"  +0.01%  `""
$scratch.Copy()
$ws.Range("E8").PasteSpecial(-4163)

$scratch.Formula = "=`"  -1.15%  `""
$scratch.Copy()
$ws.Range("E9").PasteSpecial(-4163)

$scratch.Formula = "=`"35.25`""
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)

$scratch.Formula = "=`"  +3.38%  `""
$scratch.Copy()
$ws.Range("E10").PasteSpecial(-4163)

$scratch.Formula = "=`"0.0797`""
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)

$scratch.Formula = "=`"  +1.00%  `""
$scratch.Copy()
$ws.Range("E11").PasteSpecial(-4163)

$scratch.Formula = "=`"  +2.80%  `""
$scratch.Copy()
$ws.Range("E12").PasteSpecial(-4163)

$scratch.Formula = "=`"  +0.30%  `""
$scratch.Copy()
$ws.Range("E13").PasteSpecial(-4163)

$scratch.Formula = "=`"  +2.31%  `""
$scratch.Copy()
$ws.Range("E14").PasteSpecial(-4163)

$scratch.Formula = "=`"2.787.81`""
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)

$scratch.Formula = "=`"  +2.14%  `""
$scratch.Copy()
$ws.Range("E15").PasteSpecial(-4163)

$scratch.Formula = "=`"2.424.99`""
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)

$scratch.Formula = "=`"  +3.27%  `""
$scratch.Copy()
$ws.Range("E16").PasteSpecial(-4163)

$scratch.Formula = "=`"0.828`""
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)

$scratch.Formula = "=`"  +3.80%  `""
$scratch.Copy()
$ws.Range("E17").PasteSpecial(-4163)

$scratch.Formula = "=`"43.645.18`""
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)

$scratch.Formula = "=`"  +1.21%  `""
$scratch.Copy()
$ws.Range("E18").PasteSpecial(-4163)

$scratch.Formula = "=`"  +2.07%  `""
$scratch.Copy()
$ws.Range("E19").PasteSpecial(-4163)

$scratch.Formula = "=`"12.09`""
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)

$scratch.Formula = "=`"  -1.87%  `""
$scratch.Copy()
$ws.Range("E20").PasteSpecial(-4163)

$scratch.Formula = "=`"0.0₃0900`""
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)

$scratch.Formula = "=`"  +1.37%  `""
$scratch.Copy()
$ws.Range("E21").PasteSpecial(-4163)

$scratch.Formula = "=`"68.26`""
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)

$scratch.Formula = "=`"  +0.17%  `""
$scratch.Copy()
$ws.Range("E22").PasteSpecial(-4163)

$scratch.Formula = "=`"238.28`""
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)

$scratch.Formula = "=`"  +1.16%  `""
$scratch.Copy()
$ws.Range("E23").PasteSpecial(-4163)

$scratch.Formula = "=`"  +1.07%  `""
$scratch.Copy()
$ws.Range("E24").PasteSpecial(-4163)

$scratch.Formula = "=`"  +0.87%  `""
$scratch.Copy()
$ws.Range("E25").PasteSpecial(-4163)

$scratch.Formula = "=`"  +0.25%  `""
$scratch.Copy()
$ws.Range("E26").PasteSpecial(-4163)

$scratch.Formula = "=`"24.97`""
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)

$scratch.Formula = "=`"  +1.29%  `""
$scratch.Copy()
$ws.Range("E27").PasteSpecial(-4163)

$scratch.Formula = "=`"  -0.76%  `""
$scratch.Copy()
$ws.Range("E28").PasteSpecial(-4163)

$scratch.Formula = "=`"9.42`""
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)

$scratch.Formula = "=`"  +3.14%  `""
$scratch.Copy()
$ws.Range("E29").PasteSpecial(-4163)

$scratch.Formula = "=`"32.34`""
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)

$scratch.Formula = "=`"  +3.03%  `""
$scratch.Copy()
$ws.Range("E30").PasteSpecial(-4163)

$scratch.Formula = "=`"0.119`""
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163)

$scratch.Formula = "=`"  +18.00%  `""
$scratch.Copy()
$ws.Range("E31").PasteSpecial(-4163)

$scratch.Formula = "=`"5.15`""
$scratch.Copy()
$ws.Range("D32").PasteSpecial(-4163)

$scratch.Formula = "=`"  +1.18%  `""
$scratch.Copy()
$ws.Range("E32").PasteSpecial(-4163)

$scratch.Formula = "=`"18.39`""
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)

$scratch.Formula = "=`"  +6.86%  `""
$scratch.Copy()
$ws.Range("E33").PasteSpecial(-4163)

$scratch.Formula = "=`"0.0751`""
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)

$scratch.Formula = "=`"  +3.08%  `""
$scratch.Copy()
$ws.Range("E35").PasteSpecial(-4163)

$scratch.Formula = "=`"  +3.51%  `""
$scratch.Copy()
$ws.Range("E36").PasteSpecial(-4163)

$scratch.Formula = "=`"129.98`""
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)

$scratch.Formula = "=`"  +26.38%  `""
$scratch.Copy()
$ws.Range("E37").PasteSpecial(-4163)

$scratch.Formula = "=`"2.90`""
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)

$scratch.Formula = "=`"  +4.75%  `""
$scratch.Copy()
$ws.Range("E38").PasteSpecial(-4163)

$scratch.Formula = "=`"4.39`""
$scratch.Copy()
$ws.Range("D39").PasteSpecial(-4163)

$scratch.Formula = "=`"  +0.21%  `""
$scratch.Copy()
$ws.Range("E39").PasteSpecial(-4163)

$scratch.Formula = "=`"  -1.11%  `""
$scratch.Copy()
$ws.Range("E40").PasteSpecial(-4163)

$scratch.Formula = "=`"  +0.23%  `""
$scratch.Copy()
$ws.Range("E41").PasteSpecial(-4163)

$scratch.Formula = "=`"21.23`""
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)

$scratch.Formula = "=`"  -4.73%  `""
$scratch.Copy()
$ws.Range("E42").PasteSpecial(-4163)

$scratch.Formula = "=`"1.948.65`""
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)

$scratch.Formula = "=`"  +0.38%  `""
$scratch.Copy()
$ws.Range("E43").PasteSpecial(-4163)

$scratch.Formula = "=`"  +1.31%  `""
$scratch.Copy()
$ws.Range("E44").PasteSpecial(-4163)

$scratch.Formula = "=`"  +2.26%  `""
$scratch.Copy()
$ws.Range("E45").PasteSpecial(-4163)

$scratch.Formula = "=`"  +3.44%  `""
$scratch.Copy()
$ws.Range("E46").PasteSpecial(-4163)

$scratch.Formula = "=`"9.29`""
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)

$scratch.Formula = "=`"  -1.74%  `""
$scratch.Copy()
$ws.Range("E47").PasteSpecial(-4163)

$scratch.Formula = "=`"2.642.10`""
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)

$scratch.Formula = "=`"  +2.03%  `""
$scratch.Copy()
$ws.Range("E48").PasteSpecial(-4163)

$scratch.Formula = "=`"1.56`""
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)

$scratch.Formula = "=`"  +3.65%  `""
$scratch.Copy()
$ws.Range("E49").PasteSpecial(-4163)

$scratch.Formula = "=`"52.73`""
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)

$scratch.Formula = "=`"  -0.25%  `""
$scratch.Copy()
$ws.Range("E50").PasteSpecial(-4163)

$scratch.Formula = "=`"72.38`""
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)

$scratch.Formula = "=`"  +0.21%  `""
$scratch.Copy()
$ws.Range("E51").PasteSpecial(-4163)

$scratch.Clear()
$excel.CutCopyMode = 0
